# chore: Update contact page layout and content
#
# Adds a new "2024-08-21" sheet (archived tasks) after "2024-08-20",
# appends a new archived task to the "2024-08-20" sheet, and populates
# the new sheet with its own archived tasks, matching the style used
# by the header cell ("Archived Tasks") on the other archived-task sheets.

$wb = $excel.ActiveWorkbook

# --- Update existing "2024-08-20" sheet: add row 3 ---
$sheet20 = $wb.Worksheets.Item("2024-08-20")
$sheet20.Range("A3").Value = "Portfolio Website Done ✅"

# --- Add the new "2024-08-21" sheet right after "2024-08-20" ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet20)
$newSheet.Name = "2024-08-21"

# Copy the header cell's formatting (style index 2) from the "2024-08-20"
# sheet onto the new sheet's A1, then set the cell values.
$sheet20.Range("A1").Copy() | Out-Null
$newSheet.Range("A1").PasteSpecial(-4122) | Out-Null

$newSheet.Range("A1").Value = "Archived Tasks"
$newSheet.Range("A2").Value = "Update Upwork Profile"
$newSheet.Range("A3").Value = "Fix Upwork To Show For Every Body"
